# fixed harvester column in rnasamples -- holly added S.GISH to harvester in bioSamples
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "harvester" values (column B, rows 2-25) with the new value "S.GISH"
$ws.Range("B2:B25").Value = "S.GISH"

# Cosmetic side-effects of editing that column in Excel: the column narrows/widens
# to fit, and the selection ends up on column B.
$ws.Columns("B:B").ColumnWidth = 8
$ws.Columns("B:B").Select() | Out-Null
